$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy (not move) the localization data that lives on the per-tournament
# sheets into the shared "tournament" table, so it can be looked up by key.
# New rows are appended to the bottom of the existing table data first
# (venue-key.*), then the competition-key / host-key rows are inserted at
# the top (this insertion order matters: it is what determines the order
# new strings land in the shared-strings table, matching how this was
# actually authored in Excel).

$venueKeys = @(
    @("venue-key.1", "ru-moscow_luzhniki"),
    @("venue-key.2", "ru-ekaterinburg"),
    @("venue-key.3", "ru-saint-petersburg"),
    @("venue-key.4", "ru-sochi"),
    @("venue-key.5", "ru-kazan"),
    @("venue-key.6", "ru-moscow_otkrytiye"),
    @("venue-key.7", "ru-saransk"),
    @("venue-key.8", "ru-kaliningrad"),
    @("venue-key.9", "ru-samara"),
    @("venue-key.10", "ru-rostov-on-don"),
    @("venue-key.11", "ru-nizhny-novgorod"),
    @("venue-key.12", "ru-volgograd")
)

$startRow = 17
for ($i = 0; $i -lt $venueKeys.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $venueKeys[$i][0]
    $ws.Cells.Item($r, 2).Value = $venueKeys[$i][1]
}

# Insert two new rows after the header for competition-key / host-key.
$ws.Rows("2:3").Insert()
$ws.Range("A2").Value = "competition-key"
$ws.Range("B2").Value = "mens-world-cup"
$ws.Range("A3").Value = "host-key"
$ws.Range("B3").Value = "russia"

# Grow the "tournament" table to cover the newly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I30"))

# Match the reported selection after the edit.
$null = $ws.Range("A2:XFD3").Select()
